$wb = $excel.ActiveWorkbook

# --- Sheet: Input Parameters ---
$ws1 = $wb.Worksheets.Item("Input Parameters")
$ws1.Range("B2").Value = 31.68056
$ws1.Range("C2").Value = -7.59583
$ws1.Range("B3").Value = 31.68056
$ws1.Range("C3").Value = -7.59583
$ws1.Range("B4").Value = 31.68056
$ws1.Range("C4").Value = -7.59583
$ws1.Range("B5").Value = 31.68056
$ws1.Range("C5").Value = -7.59583
$ws1.Range("C6").Value = -7.59583
$ws1.Range("C7").Value = -7.59583
$ws1.Range("C9").Value = -7.59583
$ws1.Range("B10").Value = 31.68056
$ws1.Range("C10").Value = -7.59583

# --- Sheet: Output Results ---
$ws2 = $wb.Worksheets.Item("Output Results")
$ws2.Range("D2").Value = 39234
$ws2.Range("E2").Value = 515
$ws2.Range("F2").Value = 1.52859279853724
$ws2.Range("D3").Value = 39590
$ws2.Range("E3").Value = 871
$ws2.Range("F3").Value = 1.027367046393602
$ws2.Range("F4").Value = 1.568234878753612
$ws2.Range("D5").Value = 40333
$ws2.Range("E5").Value = 1614
$ws2.Range("F5").Value = 2.153246342067261
$ws2.Range("F6").Value = 0.4341064463304053
$ws2.Range("F7").Value = 0.6792629133508745
$ws2.Range("F8").Value = 0.72784282473241
$ws2.Range("F9").Value = 0.4597446247869066
$ws2.Range("F10").Value = 0.9020416044696654
$ws2.Range("F11").Value = 0.1831497061510858
$ws2.Range("F12").Value = 1.691118173936799
$ws2.Range("F13").Value = 0.3280135870771285
$ws2.Range("F14").Value = 0.1901580738977416
$ws2.Range("F15").Value = 0.4056849792480352
$ws2.Range("F16").Value = 4.460284187979638
$ws2.Range("F17").Value = 3.527925250845044
$ws2.Range("F18").Value = 5.458218752051685
$ws2.Range("F19").Value = 6.36593828368605
$ws2.Range("F20").Value = 4.501058327880293
$ws2.Range("F21").Value = 3.545589493487691
$ws2.Range("F22").Value = 5.270495471766091
$ws2.Range("F23").Value = 6.147844923512341
$ws2.Range("F24").Value = 2.179977272884304
$ws2.Range("F25").Value = 1.284153232224922
$ws2.Range("F26").Value = 2.325970261920851
$ws2.Range("F27").Value = 2.685650868725979
$ws2.Range("F32").Value = 1.575718291759458
$ws2.Range("F33").Value = 0.807216653698976
$ws2.Range("F34").Value = 1.796068706048071
$ws2.Range("F35").Value = 1.209034815115379
$ws2.Range("D36").Value = 39235
$ws2.Range("E36").Value = 516
$ws2.Range("F36").Value = 1.353213898001245
$ws2.Range("D37").Value = 39590
$ws2.Range("E37").Value = 871
$ws2.Range("F37").Value = 0.7013404472768289
$ws2.Range("D38").Value = 39965
$ws2.Range("E38").Value = 1246
$ws2.Range("F38").Value = 1.603871174019468
$ws2.Range("D39").Value = 40322
$ws2.Range("E39").Value = 1603
$ws2.Range("F39").Value = 0.8425797837352049
